$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full dataset for rows 2-31 (9 new rows inserted at the top of the
# existing data, and 1 new row appended at the bottom, per the May 9th
# data refresh).
$data = @(
        @(-2.822264432907104,3.932364892959594,0.4758871763944626),
        @(-2.760175883769989,3.927277326583862,0.4960161261260509),
        @(-2.64756965637207,3.941246557235718,0.6029229372739793),
        @(-2.628955054283141,3.869387286901474,0.6939661800861359),
        @(-2.694027137756347,3.838324213027954,0.7007610917091369),
        @(-2.807751727104188,3.847099477052689,0.6076438263058663),
        @(-2.889829158782959,3.764959990978241,0.4738430827856064),
        @(-3.067947173118591,3.764124345779419,0.4743617072701454),
        @(-3.113880395889281,3.648524475097656,0.3546387374401094),
        @(-2.830316853523255,3.761647629737854,0.4785371914505959),
        @(-2.873888826370239,3.739926481246949,0.2519635170698165),
        @(-3.17759734392166,3.094935894012452,-1.055553257465361),
        @(-3.781754684448242,2.629349136352539,-1.783702611923218),
        @(-3.911327278614044,1.632070899009704,-2.105818438529968),
        @(-5.623295831680297,0.7180684566497808,-2.688220548629761),
        @(-6.422207009792328,0.5463402450084687,-3.518551957607269),
        @(-5.178159475326543,1.582323312759396,-3.006135582923893),
        @(-3.926182019710541,3.611184996366497,0.6073530614375975),
        @(-4.239339423179628,4.14080636501312,3.278265589475629),
        @(-1.282989490032158,2.413680851459504,0.6543272763490596),
        @(7.509865951538076,2.822640895843509,-1.422052669525146),
        @(2.871775150299069,4.098878204822539,-1.140487685799597),
        @(3.782487344741822,4.224413537979126,-1.672759181261064),
        @(-0.05776283740998212,3.633949923515318,0.03897095024586128),
        @(3.522613048553472,5.809941577911379,-2.533208680152897),
        @(4.382918024063111,3.592077732086182,-2.185777962207795),
        @(2.164362668991095,3.682830810546873,-3.28418397903442),
        @(1.118086504936218,3.899120330810547,-2.249887198209762),
        @(0.6892168045043942,4.255544376373291,-1.803213787078858),
        @(0.8177140951156612,4.59937185049057,-1.870238688588142)
    )

$numRows = $data.Count
$numCols = 3

for ($i = 0; $i -lt $numRows; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $numCols; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $data[$i][$j]
    }
}
